$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values that look numeric need a leading apostrophe so Excel
# keeps storing them as text (matching the source inlineStr), exactly as
# it would if a user typed them into a Text-formatted / already-text cell.
$ws.Range("D2").Value = "26.708.65"
$ws.Range("E2").Value = "  +0.42%  "
$ws.Range("D3").Value = "1.598.95"
$ws.Range("E3").Value = "  +0.93%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'211.34"
$ws.Range("E5").Value = "  +0.26%  "
$ws.Range("E6").Value = "  +1.36%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "'0.0617"
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("D10").Value = "'19.59"
$ws.Range("E10").Value = "  +0.40%  "
$ws.Range("D11").Value = "'0.0843"
$ws.Range("E11").Value = "  +1.18%  "
$ws.Range("D12").Value = "1.823.52"
$ws.Range("E12").Value = "  +1.00%  "
$ws.Range("D13").Value = "1.601.91"
$ws.Range("E13").Value = "  +1.00%  "
$ws.Range("E15").Value = "  -1.16%  "
$ws.Range("D16").Value = "'64.89"
$ws.Range("E16").Value = "  +0.62%  "
$ws.Range("D17").Value = "26.693.97"
$ws.Range("E17").Value = "  +0.31%  "
$ws.Range("D18").Value = "0.0₃0729"
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("D19").Value = "'208.58"
$ws.Range("E19").Value = "  +0.52%  "
$ws.Range("E20").Value = "  +0.03%  "
$ws.Range("E21").Value = "  +0.63%  "
$ws.Range("D22").Value = "'4.24"
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("D23").Value = "'2.34"
$ws.Range("E23").Value = "  -1.00%  "
$ws.Range("E24").Value = "  -0.23%  "
$ws.Range("D25").Value = "'145.55"
$ws.Range("E25").Value = "  -0.68%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").Value = "'7.23"
$ws.Range("E27").Value = "  -2.26%  "
$ws.Range("E28").Value = "  +1.72%  "
$ws.Range("E29").Value = "  -0.18%  "
$ws.Range("D30").Value = "'0.0506"
$ws.Range("E30").Value = "  +0.75%  "
$ws.Range("E31").Value = "  +0.28%  "
$ws.Range("E32").Value = "  -0.65%  "
$ws.Range("D33").Value = "'0.662"
$ws.Range("E33").Value = "  -2.50%  "
$ws.Range("E34").Value = "  +0.52%  "
$ws.Range("D35").Value = "1.286.16"
$ws.Range("E35").Value = "  -2.66%  "
$ws.Range("E36").Value = "  +1.36%  "
$ws.Range("E37").Value = "  -0.69%  "
$ws.Range("E38").Value = "  -0.31%  "
$ws.Range("E39").Value = "  +2.48%  "
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("D41").Value = "'5.42"
$ws.Range("E41").Value = "  +1.60%  "
$ws.Range("E42").Value = "  +1.18%  "
$ws.Range("D43").Value = "'0.787"
$ws.Range("E43").Value = "  +0.56%  "
$ws.Range("D44").Value = "'63.53"
$ws.Range("E44").Value = "  +0.20%  "
$ws.Range("D45").Value = "1.735.86"
$ws.Range("E45").Value = "  +0.99%  "
$ws.Range("D46").Value = "'0.906"
$ws.Range("E46").Value = "  +9.19%  "
$ws.Range("E47").Value = "  +0.55%  "
$ws.Range("E48").Value = "  -0.46%  "
$ws.Range("E49").Value = "  +2.61%  "
$ws.Range("E51").Value = "  +0.11%  "
